# consent.xlsx - "version 2" text changes: payment / EEG compensation wording,
# researcher contact line, registration consent statement and the
# "next three year(s)" plural fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (p4): English reward-estimate column text is unchanged, but the
# payment paragraph used to live here before; only row height shrinks a touch
# because the long payment text moved down to row 15 (see below).
$ws.Rows.Item(9).RowHeight = 225

# --- Row 15 (p8 "your reward" payment paragraph) ---
$ws.Cells.Item(15, 2).Value = "you will be compensated at the rate of 3.50 EUR per 30 min to a maximum of 80 EUR for online experiments. If you participate in both online and brain (EEG) experiments then the maximum compensation is 170 EUR at the end of the three years. You can annually receive your compensation, however you must provide us with your bank details, full name and address and live in the Netherlands."
$ws.Cells.Item(15, 3).Value = "ontvangt u een vergoeding van 3,50 euro per 30 minuten tot een maximum van 80 euro voor online onderzoeken. Als u zowel deelneemt aan online onderzoeken als aan registratie van hersenactiviteit (EEG), dan is de maximale vergoeding 170 euro in drie jaar. U kunt uw vergoeding jaarlijks ontvangen, daarvoor hebben wij echter wel uw bankgegevens, volledige naam en adres nodig. Tevens moet u in Nederland wonen. `n"
$ws.Cells.Item(15, 3).WrapText = $true
$ws.Rows.Item(15).RowHeight = 270.75

# --- Row 16 (h17 researcher contact details) - Dutch label loses a stray
# narrow no-break space before the trailing space.
$ws.Cells.Item(16, 3).Value = "Naam en contactgegevens van de onderzoeker "

# --- Row 18 (p10 consent-statement intro) ---
$ws.Cells.Item(18, 2).Value = "Please read and provide your signature or check mark to the following statement during your registration only if you approve."
$ws.Cells.Item(18, 3).Value = "Lees alstublieft de volgende verklaring en plaats alleen uw handtekening of een vinkje bij uw registratie als u deze goedkeurt."

# --- Row 26 (li7_b) - pluralise "next three year" -> "next three years"
$ws.Cells.Item(26, 2).Value = "next three years"

# --- Restore view state (matches the saved workbook's sheetView) ---
$ws.Range("B26").Select()
